# Generate Report for Handoff
# Updates the localization-status workbook to reflect that file "b.md" has been
# re-handed-off (new handoff package generated) instead of already being in sync.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet - row 3 corresponds to file "b.md"
# ---------------------------------------------------------------------------
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-02-21 03:03:29"

# ---------------------------------------------------------------------------
# zh-cn sheet - row 3 corresponds to file "b.md"
# ---------------------------------------------------------------------------
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-02-21 03:03:13"
$wsZhCn.Range("O3").Value = "'True"
$wsZhCn.Range("O3").Style = "Normal"
$wsZhCn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7d7fa767d677dbf045c08e7d8e5edf9f52164586/e2e/b.md."

# Widen the "Error Detail" column (R) so the new message is readable.
# NOTE: the ColumnWidth property is expressed in "characters" and gets an
# internal +5/6 padding added when stored as the OOXML <col width="..">
# attribute, so we dial in 39.1666... to land on a stored width of exactly 40.
$wsZhCn.Range("R1").ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# de-de sheet - row 3 corresponds to file "b.md"
# ---------------------------------------------------------------------------
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-02-21 03:03:29"
$wsDeDe.Range("O3").Value = "'True"
$wsDeDe.Range("O3").Style = "Normal"
$wsDeDe.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5b8f33dc302b5a2aa99f42855abaa4d3b6b8492e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7d7fa767d677dbf045c08e7d8e5edf9f52164586/e2e/b.md."

# Widen the "Error Detail" column (R) so the new message is readable.
# NOTE: the ColumnWidth property is expressed in "characters" and gets an
# internal +5/6 padding added when stored as the OOXML <col width="..">
# attribute, so we dial in 39.1666... to land on a stored width of exactly 40.
$wsDeDe.Range("R1").ColumnWidth = 39.16666666666667
